# Generate Report for Handoff
# The localization status report moved from "In Translation" to
# "Ready for handoff" for the c9a46d89 file in both zh-cn and de-de, and
# the handoff/generation timestamps were refreshed accordingly.

$wb = $excel.ActiveWorkbook

$statusNew = "Ready for handoff"

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew            # zh-cn status
$wsOverview.Range("F2").Value = $statusNew            # de-de status
$wsOverview.Range("G2").Value = "2016-08-19 02:55:31" # Latest HO Xliff Generate Date

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusNew                # Status
$wsZhCn.Range("H2").Value = "2016-08-19 02:55:27"     # Latest Handoff Datetime

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusNew                # Status
$wsDeDe.Range("H2").Value = "2016-08-19 02:55:31"     # Latest Handoff Datetime

# --- widen the Status columns to fit the longer "Ready for handoff" text ---
# (COM ColumnWidth is quantized to the nearest 1/6 of a character by this
#  host, so we target the closest representable width to the authored one.)
$newStatusWidth = 16 + 1/3

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusWidth  # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusWidth  # column F (de-de)
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusWidth      # column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusWidth      # column C (Status)
